$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# Row 6
$ws.Range("C6").Formula = "= 8 * 20127.75"

# Row 7
$ws.Range("B7").Formula = "=4 * 194176.23"
$ws.Range("C7").Formula = "=4 * 37921.61"

# Row 8
$ws.Range("B8").Formula = "=4*125486.01"
$ws.Range("C8").Formula = "=4 * 38627.17"

# Row 9
$ws.Range("B9").Formula = "=2*125605.79"
$ws.Range("C9").Formula = "=2 * 39277.56"

# Row 10
$ws.Range("B10").Value = 194494.19
$ws.Range("C10").Value = 38212.33

# Row 11
$ws.Range("B11").Value = 156932.69
$ws.Range("C11").Value = 37754.9

# Row 12
$ws.Range("B12").Value = 126163.42
$ws.Range("C12").Value = 39294.65

# Row 13
$ws.Range("B13").Value = 9483.93
$ws.Range("C13").Value = 10493.7746

# Row 14
$ws.Range("C14").Value = 5427.03
$ws.Range("E14").Value = 3703.44

# Row 15
$ws.Range("B15").Value = 4774.73
$ws.Range("C15").Value = 2509.35
$ws.Range("D15").Value = 14520.36
$ws.Range("E15").Value = 3334.15

# Row 16
$ws.Range("B16").Value = 1075.89
$ws.Range("C16").Value = 2508.82
$ws.Range("D16").Value = 8374.87
$ws.Range("E16").Value = 3108.3

# Update the active cell selection
$ws.Range("D6").Select()
